$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1043.3077
$ws.Range("I17").Value = 465.5
$ws.Range("K17").Value = 1396.5
$ws.Range("M17").Value = -1228.5

$ws.Range("H51").Value = 2268.75
$ws.Range("J51").Value = 1883.3334
$ws.Range("L51").Value = 1883.3334
$ws.Range("N51").Value = -2851.3334

$ws.Range("H98").Value = 11125.9375
$ws.Range("J98").Value = 0
$ws.Range("L98").Value = 0
$ws.Range("N98").ClearContents()

$ws.Range("H103").Value = 983.38464
$ws.Range("I103").Value = 862.7143
$ws.Range("J103").Value = 1124.1666
$ws.Range("K103").Value = 2588.1429
$ws.Range("L103").Value = 3372.4998
$ws.Range("M103").Value = -2002.1429
$ws.Range("N103").Value = -4544.4998

$ws.Range("H122").Value = 11125.9375
$ws.Range("J122").Value = 0
$ws.Range("L122").Value = 0
$ws.Range("N122").ClearContents()

$ws.Range("H137").Value = 3166
$ws.Range("I137").Value = 2499.25
$ws.Range("J137").Value = 4499.5
$ws.Range("K137").Value = 7497.75
$ws.Range("L137").Value = 13498.5
$ws.Range("M137").Value = -4947.75
$ws.Range("N137").Value = -18598.5

$ws.Range("H138").Value = 5357.396
$ws.Range("I138").Value = 1095.0714
$ws.Range("K138").Value = 3285.2142
$ws.Range("M138").Value = 1854.7858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5650.095
$ws.Range("I61").Value = 5393.2354
$ws.Range("K61").Value = 5393.2354
$ws.Range("M61").Value = -5181.2354

$ws.Range("H74").Value = 1663.8096
$ws.Range("I74").Value = 1681.1052
$ws.Range("K74").Value = 1681.1052
$ws.Range("M74").Value = -807.1052

$ws.Range("H77").Value = 1663.8096
$ws.Range("I77").Value = 1681.1052
$ws.Range("K77").Value = 8405.526
$ws.Range("M77").Value = -4037.526

$ws.Range("H136").Value = 5650.095
$ws.Range("I136").Value = 5393.2354
$ws.Range("K136").Value = 16179.7062
$ws.Range("M136").Value = -13629.7062

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H53").Value = 93569.664
$ws.Range("I53").Value = 80709
$ws.Range("K53").Value = 80709
$ws.Range("M53").Value = -80135

$ws.Range("H62").Value = 50000.5
$ws.Range("I62").Value = 0
$ws.Range("J62").Value = 50000.5
$ws.Range("K62").Value = 0
$ws.Range("L62").Value = 50000.5
$ws.Range("M62").ClearContents()
$ws.Range("N62").Value = -51372.5

$ws.Range("H63").Value = 50001
$ws.Range("J63").Value = 50001
$ws.Range("L63").Value = 50001
$ws.Range("N63").Value = -51373

$ws.Range("H65").Value = 50000.5
$ws.Range("I65").Value = 0
$ws.Range("J65").Value = 50000.5
$ws.Range("K65").Value = 0
$ws.Range("L65").Value = 150001.5
$ws.Range("M65").ClearContents()
$ws.Range("N65").Value = -156865.5

$ws.Range("H66").Value = 50001
$ws.Range("J66").Value = 50001
$ws.Range("L66").Value = 150003
$ws.Range("N66").Value = -156867

$ws.Range("H87").Value = 0
$ws.Range("J87").Value = 0
$ws.Range("L87").Value = 0
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 0
$ws.Range("J90").Value = 0
$ws.Range("L90").Value = 0
$ws.Range("N90").ClearContents()

$ws.Range("H118").Value = 77709
$ws.Range("J118").Value = 77709
$ws.Range("L118").Value = 77709
$ws.Range("N118").Value = -81023

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 1785.9546
$ws.Range("I107").Value = 2956.3
$ws.Range("J107").Value = 810.6667
$ws.Range("K107").Value = 8868.900000000001
$ws.Range("L107").Value = 2432.0001
$ws.Range("M107").Value = -6948.900000000001
$ws.Range("N107").Value = -6272.0001

$ws.Range("H117").Value = 1657
$ws.Range("I117").Value = 1314.5
$ws.Range("J117").Value = 1999.5
$ws.Range("K117").Value = 3943.5
$ws.Range("L117").Value = 5998.5
$ws.Range("M117").Value = -501.5
$ws.Range("N117").Value = -12882.5

$ws.Range("H132").Value = 2060.8235
$ws.Range("I132").Value = 1339.1666
$ws.Range("K132").Value = 12052.4994
$ws.Range("M132").Value = -9522.499400000001

$ws.Range("H134").Value = 2745.2222
$ws.Range("I134").Value = 2338.375
$ws.Range("K134").Value = 7015.125
$ws.Range("M134").Value = -1945.125

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 64.25
$ws.Range("I2").Value = 66.28570999999999
$ws.Range("J2").Value = 50
$ws.Range("K2").Value = 66.28570999999999
$ws.Range("L2").Value = 50
$ws.Range("M2").Value = 46.71429000000001
$ws.Range("N2").Value = -276

$ws.Range("H110").Value = 75000
$ws.Range("J110").Value = 75000
$ws.Range("L110").Value = 75000
$ws.Range("N110").Value = -83180

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H56").Value = 5000
$ws.Range("I56").Value = 5000
$ws.Range("K56").Value = 5000
$ws.Range("M56").Value = -4309

$ws.Range("H122").Value = 4511.857
$ws.Range("I122").Value = 4511.857
$ws.Range("K122").Value = 13535.571
$ws.Range("M122").Value = -11085.571

$ws.Range("H136").Value = 4217.3228
$ws.Range("I136").Value = 3148.6667
$ws.Range("K136").Value = 9446.000100000001
$ws.Range("M136").Value = -6896.000100000001

$ws.Range("H138").Value = 49646.5
$ws.Range("J138").Value = 49646.5
$ws.Range("L138").Value = 49646.5
$ws.Range("N138").Value = -59926.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 0
$ws.Range("J15").Value = 0
$ws.Range("L15").Value = 0
$ws.Range("N15").ClearContents()

$ws.Range("H28").Value = 15500
$ws.Range("J28").Value = 15500
$ws.Range("L28").Value = 15500
$ws.Range("N28").Value = -16196

$ws.Range("H54").Value = 25000
$ws.Range("I54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("M54").ClearContents()

$ws.Range("H75").Value = 37500
$ws.Range("J75").Value = 37500
$ws.Range("L75").Value = 37500
$ws.Range("N75").Value = -39372

$ws.Range("H78").Value = 37500
$ws.Range("J78").Value = 37500
$ws.Range("L78").Value = 112500
$ws.Range("N78").Value = -121860

$ws.Range("H81").Value = 2039.8
$ws.Range("I81").Value = 2122
$ws.Range("K81").Value = 4244
$ws.Range("M81").Value = -3183

$ws.Range("H84").Value = 2039.8
$ws.Range("I84").Value = 2122
$ws.Range("K84").Value = 21220
$ws.Range("M84").Value = -15916

$ws.Range("H122").Value = 7500.3
$ws.Range("I122").Value = 1667.6666
$ws.Range("K122").Value = 5002.9998
$ws.Range("M122").Value = -2552.9998

$ws.Range("H136").Value = 2561.5881
$ws.Range("I136").Value = 2389.5386
$ws.Range("K136").Value = 7168.6158
$ws.Range("M136").Value = -4618.6158
